$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.363.07"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.293.31"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'537.51"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'131.89"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "2.288.57"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'5.50"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'23.72"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "2.701.29"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "58.259.50"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "2.294.88"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "'10.58"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").Value = "'315.03"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'6.57"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'63.27"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'7.99"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "'171.02"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").Value = "'5.85"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").Value = "'0.383"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'17.88"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'3.94"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "'290.64"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").Value = "'140.38"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "'0.0953"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "'0.555"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'18.31"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'4.62"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  +0.99%  "
